$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "not found"
$ws.Range("C6").Value = "not found"
